$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 2132.0833
$ws.Range("J40").Value = 2216.1667
$ws.Range("K40").Value = 2132.0833
$ws.Range("L40").Value = 2216.1667
$ws.Range("M40").Value = -1957.0833
$ws.Range("N40").Value = -2566.1667
$ws.Range("H53").Value = 327.86667
$ws.Range("I53").Value = 241.16667
$ws.Range("J53").Value = 385.66666
$ws.Range("K53").Value = 241.16667
$ws.Range("L53").Value = 385.66666
$ws.Range("M53").Value = 395.83333
$ws.Range("N53").Value = -1659.66666
$ws.Range("H64").Value = 5071.2856
$ws.Range("I64").Value = 4666.3335
$ws.Range("K64").Value = 4666.3335
$ws.Range("M64").Value = -4418.3335
$ws.Range("H67").Value = 5071.2856
$ws.Range("I67").Value = 4666.3335
$ws.Range("K67").Value = 4666.3335
$ws.Range("M67").Value = -3808.3335
$ws.Range("H70").Value = 15107.071
$ws.Range("I70").Value = 999.5
$ws.Range("J70").Value = 17458.334
$ws.Range("K70").Value = 2998.5
$ws.Range("L70").Value = 52375.00199999999
$ws.Range("M70").Value = -2728.5
$ws.Range("N70").Value = -52915.00199999999
$ws.Range("H73").Value = 15107.071
$ws.Range("I73").Value = 999.5
$ws.Range("J73").Value = 17458.334
$ws.Range("K73").Value = 2998.5
$ws.Range("L73").Value = 52375.00199999999
$ws.Range("M73").Value = -2062.5
$ws.Range("N73").Value = -54247.00199999999
$ws.Range("H137").Value = 923
$ws.Range("I137").Value = 863.6667
$ws.Range("J137").Value = 1101
$ws.Range("K137").Value = 2591.0001
$ws.Range("L137").Value = 3303
$ws.Range("M137").Value = -41.0001000000002
$ws.Range("N137").Value = -8403
$ws.Range("H138").Value = 3026.9285
$ws.Range("J138").Value = 3428.4
$ws.Range("L138").Value = 10285.2
$ws.Range("N138").Value = -20565.2
$ws.Range("H141").Value = 1518
$ws.Range("I141").Value = 1518
$ws.Range("K141").Value = 4554
$ws.Range("M141").Value = 626
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 164.8
$ws.Range("I25").Value = 164.8
$ws.Range("K25").Value = 164.8
$ws.Range("M25").Value = 237.2
$ws.Range("H27").Value = 1000
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1368
$ws.Range("H32").Value = 2498.8572
$ws.Range("I32").Value = 2633.7576
$ws.Range("K32").Value = 2633.7576
$ws.Range("M32").Value = -2346.7576
$ws.Range("H88").Value = 1276
$ws.Range("I88").Value = 1722.25
$ws.Range("K88").Value = 1722.25
$ws.Range("M88").Value = -1316.25
$ws.Range("H91").Value = 1276
$ws.Range("I91").Value = 1722.25
$ws.Range("K91").Value = 1722.25
$ws.Range("M91").Value = -318.25
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").Value = $null
$ws.Range("H122").Value = 2622.9375
$ws.Range("J122").Value = 4456.4287
$ws.Range("L122").Value = 13369.2861
$ws.Range("N122").Value = -18269.2861
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1997
$ws.Range("I20").Value = 1997
$ws.Range("K20").Value = 1997
$ws.Range("M20").Value = -1750
$ws.Range("H94").Value = 967.1111
$ws.Range("I94").Value = 918.1177
$ws.Range("J94").Value = 1800
$ws.Range("K94").Value = 918.1177
$ws.Range("L94").Value = 1800
$ws.Range("M94").Value = -467.1177
$ws.Range("N94").Value = -2702
$ws.Range("H134").Value = 3116.5
$ws.Range("I134").Value = 3049.75
$ws.Range("J134").Value = 3250
$ws.Range("K134").Value = 9149.25
$ws.Range("L134").Value = 9750
$ws.Range("M134").Value = -6614.25
$ws.Range("N134").Value = -14820
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 47.4375
$ws.Range("I19").Value = 48
$ws.Range("J19").Value = 46.875
$ws.Range("K19").Value = 48
$ws.Range("L19").Value = 46.875
$ws.Range("M19").Value = 122
$ws.Range("N19").Value = -386.875
$ws.Range("H24").Value = 47.4375
$ws.Range("I24").Value = 48
$ws.Range("J24").Value = 46.875
$ws.Range("K24").Value = 48
$ws.Range("L24").Value = 46.875
$ws.Range("M24").Value = 122
$ws.Range("N24").Value = -386.875
$ws.Range("H125").Value = 140276
$ws.Range("J125").Value = 140276
$ws.Range("L125").Value = 140276
$ws.Range("N125").Value = -145196
$ws.Range("H132").Value = 5416
$ws.Range("I132").Value = 9500
$ws.Range("K132").Value = 28500
$ws.Range("M132").Value = -25970
$ws.Range("H134").Value = 1299.1
$ws.Range("I134").Value = 1343.4445
$ws.Range("K134").Value = 4030.3335
$ws.Range("M134").Value = -1495.3335
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 832.5
$ws.Range("I14").Value = 832.5
$ws.Range("K14").Value = 2497.5
$ws.Range("M14").Value = -2324.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 2000
$ws.Range("K27").Value = 2000
$ws.Range("M27").Value = -1834
$ws.Range("H31").Value = 2000
$ws.Range("H37").Value = 2000
$ws.Range("H80").Value = 1983
$ws.Range("I80").Value = 1983
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1983
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -985
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 1983
$ws.Range("I83").Value = 1983
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 9915
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -4923
$ws.Range("N83").Value = $null
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("H132").Value = 1197
$ws.Range("J132").Value = 1197
$ws.Range("L132").Value = 3591
$ws.Range("N132").Value = -8651
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1039.75
$ws.Range("I16").Value = 979
$ws.Range("K16").Value = 979
$ws.Range("M16").Value = -809
$ws.Range("H22").Value = 1060.1
$ws.Range("I22").Value = 1022.3333
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 1022.3333
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -727.3333
$ws.Range("N22").Value = -1990
$ws.Range("H27").Value = 1060.1
$ws.Range("I27").Value = 1022.3333
$ws.Range("J27").Value = 1400
$ws.Range("K27").Value = 1022.3333
$ws.Range("L27").Value = 1400
$ws.Range("M27").Value = -915.3333
$ws.Range("N27").Value = -1614
$ws.Range("H61").Value = 4316.3335
$ws.Range("I61").Value = 4509.4
$ws.Range("K61").Value = 4509.4
$ws.Range("M61").Value = -4307.4
$ws.Range("H68").Value = 2511.5386
$ws.Range("I68").Value = 2513.6365
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 2513.6365
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -1764.6365
$ws.Range("N68").Value = -3998
$ws.Range("H71").Value = 2511.5386
$ws.Range("I71").Value = 2513.6365
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 12568.1825
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -8824.182500000001
$ws.Range("N71").Value = -19988
$ws.Range("H93").Value = 1594.7273
$ws.Range("I93").Value = 1334
$ws.Range("J93").Value = 1907.6
$ws.Range("K93").Value = 1334
$ws.Range("L93").Value = 1907.6
$ws.Range("M93").Value = -86
$ws.Range("N93").Value = -4403.6
$ws.Range("H100").Value = 3595
$ws.Range("I100").Value = 3335.2856
$ws.Range("J100").Value = 3958.6
$ws.Range("K100").Value = 3335.2856
$ws.Range("L100").Value = 3958.6
$ws.Range("M100").Value = -2794.2856
$ws.Range("N100").Value = -5040.6
$ws.Range("H113").Value = 4316.3335
$ws.Range("I113").Value = 4509.4
$ws.Range("K113").Value = 4509.4
$ws.Range("M113").Value = -2339.4
$ws.Range("H122").Value = 4887.0835
$ws.Range("I122").Value = 4362.7896
$ws.Range("J122").Value = 6879.4
$ws.Range("K122").Value = 13088.3688
$ws.Range("L122").Value = 20638.2
$ws.Range("M122").Value = -10638.3688
$ws.Range("N122").Value = -25538.2
$ws.Range("H130").Value = 65000
$ws.Range("J130").Value = 65000
$ws.Range("L130").Value = 65000
$ws.Range("N130").Value = -75040
$ws.Range("H132").Value = 1370.75
$ws.Range("I132").Value = 1227.8334
$ws.Range("K132").Value = 3683.5002
$ws.Range("M132").Value = -1153.5002
$ws.Range("H136").Value = 2926.5715
$ws.Range("I136").Value = 1856.4
$ws.Range("J136").Value = 5602
$ws.Range("K136").Value = 5569.200000000001
$ws.Range("L136").Value = 16806
$ws.Range("M136").Value = -3019.200000000001
$ws.Range("N136").Value = -21906
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 9000
$ws.Range("I48").Value = 9000
$ws.Range("K48").Value = 9000
$ws.Range("M48").Value = -8431
$ws.Range("H62").Value = 24426.143
$ws.Range("I62").Value = 24498
$ws.Range("J62").Value = 24397.4
$ws.Range("K62").Value = 24498
$ws.Range("L62").Value = 24397.4
$ws.Range("M62").Value = -23874
$ws.Range("N62").Value = -25645.4
$ws.Range("H65").Value = 24426.143
$ws.Range("I65").Value = 24498
$ws.Range("J65").Value = 24397.4
$ws.Range("K65").Value = 122490
$ws.Range("L65").Value = 121987
$ws.Range("M65").Value = -119370
$ws.Range("N65").Value = -128227
$ws.Range("H122").Value = 2177.889
$ws.Range("I122").Value = 1657.4286
$ws.Range("J122").Value = 3999.5
$ws.Range("K122").Value = 4972.2858
$ws.Range("L122").Value = 11998.5
$ws.Range("M122").Value = -2522.2858
$ws.Range("N122").Value = -16898.5
$ws.Range("H132").Value = 3315.2
$ws.Range("I132").Value = 3599.125
$ws.Range("J132").Value = 2179.5
$ws.Range("K132").Value = 10797.375
$ws.Range("L132").Value = 6538.5
$ws.Range("M132").Value = -8267.375
$ws.Range("N132").Value = -11598.5
$ws.Range("H136").Value = 5097.5835
$ws.Range("I136").Value = 5097.5835
$ws.Range("K136").Value = 15292.7505
$ws.Range("M136").Value = -12742.7505
